$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2, shifting everything below down
$ws.Rows.Item(2).Insert()

$ws.Range("B2:H2").Merge()
$ws.Range("B2:H2").Value = "Valor do ruído variando entre ± 2%"

$ws.Range("L6").Select()
